$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Build the two new header blocks (Iteration_1 / Iteration_2),
#    merging first (like B1:D1) and then cloning the formatting of
#    the existing "Standalone" cell (B1) onto each new header cell,
#    and the "Interval" row (B2) onto the new interval-label cells.
#    One cell at a time keeps the paste from synthesising new
#    composite border styles.
# ---------------------------------------------------------------
$ws.Range("E1:G1").Merge()
$ws.Range("H1:J1").Merge()

$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("J2").PasteSpecial(-4122)

$ws.Range("E1").Value = "Iteration_1"
$ws.Range("H1").Value = "Iteration_2"

$ws.Range("E2").Value = "2030"
$ws.Range("F2").Value = "2040"
$ws.Range("G2").Value = "2050"
$ws.Range("H2").Value = "2030"
$ws.Range("I2").Value = "2040"
$ws.Range("J2").Value = "2050"

# ---------------------------------------------------------------
# 2. New numeric data for columns E:J, rows 4-16 (default 0, then
#    the two rows that actually carry non-zero iteration results).
# ---------------------------------------------------------------
$ws.Range("E4:J16").Value = 0

$ws.Range("E6").Value = 1717666.256705075
$ws.Range("F6").Value = 1291315.402767987
$ws.Range("G6").Value = 1289650.900672362
$ws.Range("H6").Value = 1717666.256704899
$ws.Range("I6").Value = 1290907.069595588
$ws.Range("J6").Value = 1289650.900680386

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 426262.425611852
$ws.Range("G8").Value = 427926.582477398
$ws.Range("H8").Value = 0.000000001185345212715367
$ws.Range("I8").Value = 426670.6740929194
$ws.Range("J8").Value = 427926.5824773965

# ---------------------------------------------------------------
# 3. Updated values in the original B:D block (re-solved numbers).
# ---------------------------------------------------------------
$ws.Range("B4").Value = 1718092.730173202
$ws.Range("C6").Value = 1018317.089681282
$ws.Range("D6").Value = 1016974.754263447
$ws.Range("C8").Value = 699204.1168262111
$ws.Range("D8").Value = 700546.1738337319
